# Update loading_percent results (rows 2-25, columns B,C,E,F,G,H,N,O)
# for the "case with 380 kV done" re-run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 11.76875250345507
$ws.Cells.Item(2, 3).Value = 9.783263982431103
$ws.Cells.Item(2, 5).Value = 16.64426960555375
$ws.Cells.Item(2, 6).Value = 16.86991607391245
$ws.Cells.Item(2, 7).Value = 16.36700783947495
$ws.Cells.Item(2, 8).Value = 10.85676619307149
$ws.Cells.Item(2, 14).Value = 15.31887641550286
$ws.Cells.Item(2, 15).Value = 14.86898206126878
$ws.Cells.Item(3, 2).Value = 11.04976545271134
$ws.Cells.Item(3, 3).Value = 9.431048272091612
$ws.Cells.Item(3, 5).Value = 15.69130998461567
$ws.Cells.Item(3, 6).Value = 15.89584955866815
$ws.Cells.Item(3, 7).Value = 16.22074674863227
$ws.Cells.Item(3, 8).Value = 10.89943263847244
$ws.Cells.Item(3, 14).Value = 15.33975267415992
$ws.Cells.Item(3, 15).Value = 14.9101061168582
$ws.Cells.Item(4, 2).Value = 10.58164607670726
$ws.Cells.Item(4, 3).Value = 9.206720496589162
$ws.Cells.Item(4, 5).Value = 15.08038736009158
$ws.Cells.Item(4, 6).Value = 15.26997757108489
$ws.Cells.Item(4, 7).Value = 16.14055115690676
$ws.Cells.Item(4, 8).Value = 10.92813508225121
$ws.Cells.Item(4, 14).Value = 15.35435419603894
$ws.Cells.Item(4, 15).Value = 14.94063351715365
$ws.Cells.Item(5, 2).Value = 10.38420636238841
$ws.Cells.Item(5, 3).Value = 9.113373427048833
$ws.Cells.Item(5, 5).Value = 14.82520848471167
$ws.Cells.Item(5, 6).Value = 15.008197319934
$ws.Cells.Item(5, 7).Value = 16.11032450579814
$ws.Cells.Item(5, 8).Value = 10.94045961441528
$ws.Cells.Item(5, 14).Value = 15.36075373011077
$ws.Cells.Item(5, 15).Value = 14.9543925104811
$ws.Cells.Item(6, 2).Value = 10.35101932346483
$ws.Cells.Item(6, 3).Value = 9.097759585916151
$ws.Cells.Item(6, 5).Value = 14.78246917209293
$ws.Cells.Item(6, 6).Value = 14.96433081551589
$ws.Cells.Item(6, 7).Value = 16.10545453498606
$ws.Cells.Item(6, 8).Value = 10.94254397153341
$ws.Cells.Item(6, 14).Value = 15.36184352840801
$ws.Cells.Item(6, 15).Value = 14.956756586761
$ws.Cells.Item(7, 2).Value = 10.57901032984231
$ws.Cells.Item(7, 3).Value = 9.205469271009816
$ws.Cells.Item(7, 5).Value = 15.07697073574904
$ws.Cells.Item(7, 6).Value = 15.26647399323133
$ws.Cells.Item(7, 7).Value = 16.14013353104891
$ws.Cells.Item(7, 8).Value = 10.92829875474677
$ws.Cells.Item(7, 14).Value = 15.35443868212986
$ws.Cells.Item(7, 15).Value = 14.94081374705108
$ws.Cells.Item(8, 2).Value = 11.5263914581252
$ws.Cells.Item(8, 3).Value = 9.663548666584678
$ws.Cells.Item(8, 5).Value = 16.32118758913177
$ws.Cells.Item(8, 6).Value = 16.5399640634477
$ws.Cells.Item(8, 7).Value = 16.31460732282113
$ws.Cells.Item(8, 8).Value = 10.87095653797922
$ws.Cells.Item(8, 14).Value = 15.32570491134944
$ws.Cells.Item(8, 15).Value = 14.88206165448343
$ws.Cells.Item(9, 2).Value = 13.17199167523131
$ws.Cells.Item(9, 3).Value = 10.49413797716402
$ws.Cells.Item(9, 5).Value = 18.67978125205491
$ws.Cells.Item(9, 6).Value = 19.00274580682531
$ws.Cells.Item(9, 7).Value = 16.73098608743184
$ws.Cells.Item(9, 8).Value = 10.77846819091577
$ws.Cells.Item(9, 14).Value = 15.2834718158024
$ws.Cells.Item(9, 15).Value = 14.80905814384943
$ws.Cells.Item(10, 2).Value = 14.25123342279154
$ws.Cells.Item(10, 3).Value = 11.05883435247552
$ws.Cells.Item(10, 5).Value = 20.33890945944937
$ws.Cells.Item(10, 6).Value = 20.67494806633232
$ws.Cells.Item(10, 7).Value = 17.07899124896133
$ws.Cells.Item(10, 8).Value = 10.72279743622195
$ws.Cells.Item(10, 14).Value = 15.26099783240557
$ws.Cells.Item(10, 15).Value = 14.78157717454036
$ws.Cells.Item(11, 2).Value = 14.71404772981935
$ws.Cells.Item(11, 3).Value = 11.3051191318475
$ws.Cells.Item(11, 5).Value = 21.05096805927518
$ws.Cells.Item(11, 6).Value = 21.3917225636224
$ws.Cells.Item(11, 7).Value = 17.2456579699058
$ws.Cells.Item(11, 8).Value = 10.70016168518192
$ws.Cells.Item(11, 14).Value = 15.25262046682644
$ws.Cells.Item(11, 15).Value = 14.77483138689691
$ws.Cells.Item(12, 2).Value = 14.88526255942038
$ws.Cells.Item(12, 3).Value = 11.39680681745921
$ws.Cells.Item(12, 5).Value = 21.31450442040589
$ws.Cells.Item(12, 6).Value = 21.65686569030329
$ws.Cells.Item(12, 7).Value = 17.30990001103461
$ws.Cells.Item(12, 8).Value = 10.69197874622983
$ws.Cells.Item(12, 14).Value = 15.24971268020281
$ws.Cells.Item(12, 15).Value = 14.77310970222185
$ws.Cells.Item(13, 2).Value = 14.84856810845013
$ws.Cells.Item(13, 3).Value = 11.37713102108689
$ws.Cells.Item(13, 5).Value = 21.25801812906082
$ws.Cells.Item(13, 6).Value = 21.60004134736742
$ws.Cells.Item(13, 7).Value = 17.29601531870635
$ws.Cells.Item(13, 8).Value = 10.69372376772159
$ws.Cells.Item(13, 14).Value = 15.25032717381931
$ws.Cells.Item(13, 15).Value = 14.77344338470698
$ws.Cells.Item(14, 2).Value = 14.72821467350176
$ws.Cells.Item(14, 3).Value = 11.31269416584684
$ws.Cells.Item(14, 5).Value = 21.07277151160761
$ws.Cells.Item(14, 6).Value = 21.4136618050453
$ws.Cells.Item(14, 7).Value = 17.25092100311338
$ws.Cells.Item(14, 8).Value = 10.69948066764736
$ws.Cells.Item(14, 14).Value = 15.25237594516685
$ws.Cells.Item(14, 15).Value = 14.77467302544808
$ws.Cells.Item(15, 2).Value = 14.65396836194949
$ws.Cells.Item(15, 3).Value = 11.27301822421703
$ws.Cells.Item(15, 5).Value = 20.95850858006504
$ws.Cells.Item(15, 6).Value = 21.29868154950795
$ws.Cells.Item(15, 7).Value = 17.22344424352857
$ws.Cells.Item(15, 8).Value = 10.70305761877106
$ws.Cells.Item(15, 14).Value = 15.25366529938176
$ws.Cells.Item(15, 15).Value = 14.77553481126771
$ws.Cells.Item(16, 2).Value = 14.22042142469198
$ws.Cells.Item(16, 3).Value = 11.04252127669481
$ws.Cells.Item(16, 5).Value = 20.29151865693921
$ws.Cells.Item(16, 6).Value = 20.62722412089977
$ws.Cells.Item(16, 7).Value = 17.06826120938054
$ws.Cells.Item(16, 8).Value = 10.72433102283695
$ws.Cells.Item(16, 14).Value = 15.26158238313086
$ws.Cells.Item(16, 15).Value = 14.78213430153499
$ws.Cells.Item(17, 2).Value = 13.94724897119194
$ws.Cells.Item(17, 3).Value = 10.89836581732628
$ws.Cells.Item(17, 5).Value = 19.87143271024959
$ws.Cells.Item(17, 6).Value = 20.20408069597325
$ws.Cells.Item(17, 7).Value = 16.97515124039662
$ws.Cells.Item(17, 8).Value = 10.73807177962918
$ws.Cells.Item(17, 14).Value = 15.26691144763762
$ws.Cells.Item(17, 15).Value = 14.78766090471252
$ws.Cells.Item(18, 2).Value = 13.78747849111957
$ws.Cells.Item(18, 3).Value = 10.81445640698349
$ws.Cells.Item(18, 5).Value = 19.62579047592796
$ws.Cells.Item(18, 6).Value = 19.95656407809801
$ws.Cells.Item(18, 7).Value = 16.92238808757671
$ws.Cells.Item(18, 8).Value = 10.74622809576839
$ws.Cells.Item(18, 14).Value = 15.27015041192342
$ws.Cells.Item(18, 15).Value = 14.79138112824096
$ws.Cells.Item(19, 2).Value = 13.73292794641661
$ws.Cells.Item(19, 3).Value = 10.78587689919378
$ws.Cells.Item(19, 5).Value = 19.54192878574668
$ws.Cells.Item(19, 6).Value = 19.87204792380568
$ws.Cells.Item(19, 7).Value = 16.90466145110535
$ws.Cells.Item(19, 8).Value = 10.74903307423667
$ws.Cells.Item(19, 14).Value = 15.27127695263133
$ws.Cells.Item(19, 15).Value = 14.79273355951095
$ws.Cells.Item(20, 2).Value = 13.97660278242479
$ws.Cells.Item(20, 3).Value = 10.91381478267117
$ws.Cells.Item(20, 5).Value = 19.91656737790698
$ws.Cells.Item(20, 6).Value = 20.24955283636154
$ws.Cells.Item(20, 7).Value = 16.98498158251252
$ws.Cells.Item(20, 8).Value = 10.73658285360295
$ws.Cells.Item(20, 14).Value = 15.26632617581701
$ws.Cells.Item(20, 15).Value = 14.7870165085575
$ws.Cells.Item(21, 2).Value = 14.76367507730339
$ws.Cells.Item(21, 3).Value = 11.33166394944583
$ws.Cells.Item(21, 5).Value = 21.12734834900369
$ws.Cells.Item(21, 6).Value = 21.46857628470577
$ws.Cells.Item(21, 7).Value = 17.26413625239926
$ws.Cells.Item(21, 8).Value = 10.6977791605729
$ws.Cells.Item(21, 14).Value = 15.25176699943091
$ws.Cells.Item(21, 15).Value = 14.77428921083035
$ws.Cells.Item(22, 2).Value = 15.25451062510966
$ws.Cells.Item(22, 3).Value = 11.59555325538767
$ws.Cells.Item(22, 5).Value = 21.88309581376297
$ws.Cells.Item(22, 6).Value = 22.22866616901552
$ws.Cells.Item(22, 7).Value = 17.45312470743625
$ws.Cells.Item(22, 8).Value = 10.674685431394
$ws.Cells.Item(22, 14).Value = 15.24379322122525
$ws.Cells.Item(22, 15).Value = 14.77082722390246
$ws.Cells.Item(23, 2).Value = 14.99469653958081
$ws.Cells.Item(23, 3).Value = 11.45556737675912
$ws.Cells.Item(23, 5).Value = 21.48298259774992
$ws.Cells.Item(23, 6).Value = 21.82633154458858
$ws.Cells.Item(23, 7).Value = 17.35168398073238
$ws.Cells.Item(23, 8).Value = 10.68680292095143
$ws.Cells.Item(23, 14).Value = 15.24790824784766
$ws.Cells.Item(23, 15).Value = 14.77222911523444
$ws.Cells.Item(24, 2).Value = 13.96334039684631
$ws.Cells.Item(24, 3).Value = 10.90683351763757
$ws.Cells.Item(24, 5).Value = 19.89617485944633
$ws.Cells.Item(24, 6).Value = 20.22900810905287
$ws.Cells.Item(24, 7).Value = 16.98053489009585
$ws.Cells.Item(24, 8).Value = 10.73725519796138
$ws.Cells.Item(24, 14).Value = 15.26659023133348
$ws.Cells.Item(24, 15).Value = 14.78730614927234
$ws.Cells.Item(25, 2).Value = 12.74959633622283
$ws.Cells.Item(25, 3).Value = 10.27716578769649
$ws.Cells.Item(25, 5).Value = 18.03071496381154
$ws.Cells.Item(25, 6).Value = 18.34778573295695
$ws.Cells.Item(25, 7).Value = 16.61070058034226
$ws.Cells.Item(25, 8).Value = 10.80134022437451
$ws.Cells.Item(25, 14).Value = 15.29339101683061
$ws.Cells.Item(25, 15).Value = 14.82424253864333
